# The dataset-statistics sheet listed one row per mutated-bug id with the
# involving-feature set for that bug. The bug "_MultipleBugs_.NOB_1.ID_134"
# (features == ["Extract"]) is removed from the table; every row below it
# shifts up by one, and the shared "NUM OF INVOLVING FEATURES" formula /
# table dimension follow along automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12 holds "_MultipleBugs_.NOB_1.ID_134" -> ["Extract"]; delete it so the
# rows below (ID_153, ID_190, ...) move up and the sheet shrinks from
# A1:C57 to A1:C56.
$ws.Rows.Item(12).Delete()

# Reflect where the author's cursor ended up after the edit.
$ws.Range("G57").Select()
$excel.ActiveWindow.ScrollRow = 37
$excel.ActiveWindow.ScrollColumn = 1
